$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}


# Row 2 - BNB
Set-TextValue "D2" "309.09"
Set-TextValue "E2" "0.15%"

# Row 3 - OKB
Set-TextValue "D3" "41.09"
Set-TextValue "E3" "-0.93%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.243"
Set-TextValue "E4" "2.36%"

# Row 5 - Cronos
Set-TextValue "D5" "0.07662"
Set-TextValue "E5" "0.61%"

# Row 6 - GateToken (was FTXToken)
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "4.317"
Set-TextValue "E6" "1.53%"

# Row 7 - FTXToken (was MXToken)
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.624"
Set-TextValue "E7" "0.07%"

# Row 8 - BTSEToken
Set-TextValue "D8" "2.488"
Set-TextValue "E8" "-0.07%"

# Row 9 - MXToken (was LiechtensteinCryptoassetsExchange)
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D9" "0.9194"
Set-TextValue "E9" "1.63%"

# Row 10 - LiechtensteinCryptoassetsExchange (was WazirX)
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D10" "0.1246"
Set-TextValue "E10" "15.39%"

# Row 11 - WazirX (was MandalaExchangeToken)
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1840"
Set-TextValue "E11" "4.17%"

# Row 12 - MandalaExchangeToken (was BitrueCoin)
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.09170"
Set-TextValue "E12" "0.35%"

# Row 13 - BitrueCoin (was BitMartToken)
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.04292"
Set-TextValue "E13" "0.67%"

# Row 14 - BitMartToken (was BitForexToken)
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.1052"
Set-TextValue "E14" "0.17%"

# Row 15 - BitForexToken (was TigerCash)
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001261"
Set-TextValue "E15" "0.31%"

# Row 16 - TigerCash (was UpBots)
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.005743"
Set-TextValue "E16" "-1.75%"

# Row 17 - UpBots (was LEO)
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D17" "0.007498"
Set-TextValue "E17" "2,392.15%"

# Row 18 - LEO (was GateToken)
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D18" "3.354"
Set-TextValue "E18" "-0.25%"

# Row 19 - BitpandaEcosystemToken
Set-TextValue "E19" "1.21%"

# Row 20 - MCDex
Set-TextValue "D20" "7.193"
Set-TextValue "E20" "10.12%"

# Row 21 - ProBitToken
Set-TextValue "D21" "0.1384"
Set-TextValue "E21" "1.45%"

# Row 22 - ZBToken
Set-TextValue "E22" "9.20%"

# Row 23 - CoinExToken
Set-TextValue "D23" "0.04068"
Set-TextValue "E23" "-2.96%"

# Row 24 - BitKan
Set-TextValue "D24" "0.001261"
Set-TextValue "E24" "3.31%"

# Row 25 - HotbitToken
Set-TextValue "E25" "1.67%"

# Row 26 - NitroEx
Set-TextValue "D26" "0.0001273"
Set-TextValue "E26" "-1.96%"

# Row 38 - One
Set-TextValue "E38" "1.75%"

# Row 39 - IDEX
Set-TextValue "D39" "0.05309"
Set-TextValue "E39" "2.18%"

# Row 40 - KickToken
Set-TextValue "D40" "0.007844"
Set-TextValue "E40" "0.92%"

# Row 41 - BKEXToken
Set-TextValue "E41" "1.10%"

# Row 42 - Dexo
Set-TextValue "D42" "0.006823"
Set-TextValue "E42" "-1.82%"

# Row 43 - CEJI
Set-TextValue "D43" "0.001906"
Set-TextValue "E43" "-0.60%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.008370"
Set-TextValue "E44" "3.89%"

# Row 45 - PooCoin
Set-TextValue "D45" "0.3066"
Set-TextValue "E45" "0.15%"

# Row 46 - CoinLion
Set-TextValue "D46" "0.00006652"
Set-TextValue "E46" "-1.23%"

# Row 47 - Kangarootoken
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "0.34%"

# Row 48 - BOLO
Set-TextValue "E48" "1,763.72%"

# Row 50 - CryptobidCoin
Set-TextValue "D50" "0.00002105"
Set-TextValue "E50" "0.34%"

# Row 51 - SpecialPowerGold
Set-TextValue "D51" "0.0002005"
Set-TextValue "E51" "0.34%"
